$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7493.107
$ws.Range("I62").Value = 7597.9414
$ws.Range("J62").Value = 7331.091
$ws.Range("K62").Value = 7597.9414
$ws.Range("L62").Value = 7331.091
$ws.Range("M62").Value = -6973.9414
$ws.Range("N62").Value = -8579.091
$ws.Range("H65").Value = 7493.107
$ws.Range("I65").Value = 7597.9414
$ws.Range("J65").Value = 7331.091
$ws.Range("K65").Value = 37989.70699999999
$ws.Range("L65").Value = 36655.455
$ws.Range("M65").Value = -34869.70699999999
$ws.Range("N65").Value = -42895.455
$ws.Range("H134").Value = 48379.8
$ws.Range("J134").Value = 48379.8
$ws.Range("L134").Value = 48379.8
$ws.Range("N134").Value = -58519.8

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 269.4
$ws.Range("J5").Value = 289
$ws.Range("L5").Value = 289
$ws.Range("N5").Value = -513
$ws.Range("H28").Value = 1460.5
$ws.Range("I28").Value = 1460.5
$ws.Range("K28").Value = 1460.5
$ws.Range("M28").Value = -1268.5
$ws.Range("H99").Value = 1460.5
$ws.Range("I99").Value = 1460.5
$ws.Range("K99").Value = 1460.5
$ws.Range("M99").Value = 1534.5
$ws.Range("H108").Value = 20132.8
$ws.Range("J108").Value = 20132.8
$ws.Range("L108").Value = 20132.8
$ws.Range("N108").Value = -27812.8
$ws.Range("H110").Value = 1985.5
$ws.Range("I110").Value = 1469.8889
$ws.Range("J110").Value = 2648.4285
$ws.Range("K110").Value = 1469.8889
$ws.Range("L110").Value = 2648.4285
$ws.Range("M110").Value = 575.1111000000001
$ws.Range("N110").Value = -6738.4285
$ws.Range("H132").Value = 1839438.8
$ws.Range("I132").Value = 1114.2222
$ws.Range("J132").Value = 11766391
$ws.Range("K132").Value = 3342.6666
$ws.Range("L132").Value = 35299173
$ws.Range("M132").Value = -812.6665999999996
$ws.Range("N132").Value = -35304233

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 269.4
$ws.Range("J4").Value = 289
$ws.Range("L4").Value = 289
$ws.Range("N4").Value = -519
$ws.Range("H56").Value = 22995
$ws.Range("J56").Value = 22995
$ws.Range("L56").Value = 22995
$ws.Range("N56").Value = -24473
$ws.Range("H98").Value = 30000
$ws.Range("I98").Value = 30000
$ws.Range("K98").Value = 30000
$ws.Range("M98").Value = -27005
$ws.Range("H118").Value = 49800
$ws.Range("J118").Value = 49800
$ws.Range("L118").Value = 49800
$ws.Range("N118").Value = -53114

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 15788.5
$ws.Range("I33").Value = 2731
$ws.Range("J33").Value = 18400
$ws.Range("K33").Value = 2731
$ws.Range("L33").Value = 18400
$ws.Range("M33").Value = -2352
$ws.Range("N33").Value = -19158
$ws.Range("H76").Value = 2000
$ws.Range("I76").Value = 2000
$ws.Range("K76").Value = 2000
$ws.Range("M76").Value = -1685
$ws.Range("H79").Value = 2000
$ws.Range("I79").Value = 2000
$ws.Range("K79").Value = 2000
$ws.Range("M79").Value = -908
$ws.Range("H92").Value = 18800.5
$ws.Range("J92").Value = 18800.5
$ws.Range("L92").Value = 18800.5
$ws.Range("N92").Value = -23792.5
$ws.Range("H122").Value = 17857716
$ws.Range("I122").Value = 25000442
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 75001326
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -74998876
$ws.Range("N122").Value = -7600

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H105").Value = 7666.6665
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 7666.6665
$ws.Range("K105").Value = 0
$ws.Range("L105").ClearContents()
$ws.Range("M105").Value = 22999.9995
$ws.Range("N105").Value = -28241.9995
$ws.Range("H131").Value = 798.97
$ws.Range("I131").Value = 491.42856
$ws.Range("J131").Value = 822.1183
$ws.Range("K131").Value = 1474.28568
$ws.Range("L131").Value = 2466.3549
$ws.Range("M131").Value = 3565.71432
$ws.Range("N131").Value = -12546.3549

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 2760
$ws.Range("I36").Value = 2760
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2760
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -2275
$ws.Range("H43").Value = 4974.8
$ws.Range("I43").Value = 1927.5
$ws.Range("J43").Value = 7006.3335
$ws.Range("K43").Value = 1927.5
$ws.Range("L43").Value = 7006.3335
$ws.Range("M43").Value = -1776.5
$ws.Range("N43").Value = -7308.3335
$ws.Range("H46").Value = 3000
$ws.Range("J46").Value = 3000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3312
$ws.Range("H64").Value = 19800
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 19800
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H70").Value = 9863.412
$ws.Range("I70").Value = 12156.583
$ws.Range("J70").Value = 4359.8
$ws.Range("K70").Value = 12156.583
$ws.Range("L70").Value = 4359.8
$ws.Range("M70").Value = -11886.583
$ws.Range("N70").Value = -4899.8
$ws.Range("H73").Value = 9863.412
$ws.Range("I73").Value = 12156.583
$ws.Range("J73").Value = 4359.8
$ws.Range("K73").Value = 12156.583
$ws.Range("L73").Value = 4359.8
$ws.Range("M73").Value = -11220.583
$ws.Range("N73").Value = -6231.8

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1153268
$ws.Range("I22").Value = 2111651.2
$ws.Range("J22").Value = 3208
$ws.Range("K22").Value = 2111651.2
$ws.Range("L22").Value = 3208
$ws.Range("M22").Value = -2111356.2
$ws.Range("N22").Value = -3798
$ws.Range("H27").Value = 1153268
$ws.Range("I27").Value = 2111651.2
$ws.Range("J27").Value = 3208
$ws.Range("K27").Value = 2111651.2
$ws.Range("L27").Value = 3208
$ws.Range("M27").Value = -2111544.2
$ws.Range("N27").Value = -3422
$ws.Range("H54").Value = 8000
$ws.Range("J54").Value = 8000
$ws.Range("L54").Value = 8000
$ws.Range("N54").Value = -9288

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 21720.285
$ws.Range("J109").Value = 21720.285
$ws.Range("L109").Value = 21720.285
$ws.Range("N109").Value = -24494.285

Write-Output "Applied changes: 165 set, 5 cleared"